$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AW1").Value = 'Odd_CS_3-3_HT'
$ws.Range("AX1").Value = 'Odd_CS_0-1_HT'
$ws.Range("AY1").Value = 'Odd_CS_0-2_HT'
$ws.Range("AZ1").Value = 'Odd_CS_1-2_HT'
$ws.Range("BA1").Value = 'Odd_CS_0-3_HT'
$ws.Range("BB1").Value = 'Odd_CS_1-3_HT'
$ws.Range("BC1").Value = 'Odd_CS_2-3_HT'
$ws.Range("A2").Value = '2iHAeNhT'
$ws.Range("C2").Value = '05:35'
$ws.Range("E2").Value = 'Adelaide United'
$ws.Range("F2").Value = 'Western United'
$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 2.3
$ws.Range("K2").Value = 2.6
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 1.13
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 1.44
$ws.Range("R2").Value = 2.75
$ws.Range("S2").Value = 1.22
$ws.Range("T2").Value = 4
$ws.Range("U2").Value = 1.44
$ws.Range("V2").Value = 2.63
$ws.Range("W2").Value = 13
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 12
$ws.Range("AB2").Value = 17
$ws.Range("AC2").Value = 23
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 29
$ws.Range("AG2").Value = 81
$ws.Range("AI2").Value = 26
$ws.Range("AJ2").Value = 13
$ws.Range("AK2").Value = 41
$ws.Range("AL2").Value = 26
$ws.Range("AM2").Value = 23
$ws.Range("AN2").Value = 4.33
$ws.Range("AO2").Value = 9
$ws.Range("AQ2").Value = 26
$ws.Range("AS2").Value = 67
$ws.Range("AT2").Value = 4
$ws.Range("AU2").Value = 7
$ws.Range("AV2").Value = 34
$ws.Range("AW2").Value = 251
$ws.Range("AX2").Value = 6.5
$ws.Range("AY2").Value = 19
$ws.Range("AZ2").Value = 21
$ws.Range("BA2").Value = 51
$ws.Range("BB2").Value = 51
$ws.Range("BC2").Value = 101
$ws.Range("G3").Value = 2.63
$ws.Range("I3").Value = 2.75
$ws.Range("J3").Value = 3.25
$ws.Range("L3").Value = 3.4
$ws.Range("N3").Value = 10
$ws.Range("S3").Value = 1.4
$ws.Range("T3").Value = 2.75
$ws.Range("Z3").Value = 26
$ws.Range("AF3").Value = 41
$ws.Range("AH3").Value = 9
$ws.Range("AN3").Value = 4.5
$ws.Range("AT3").Value = 2.75
$ws.Range("AW3").Value = 501
$ws.Range("AX3").Value = 4.75
$ws.Range("AY3").Value = 15
$ws.Range("AZ3").Value = 23
$ws.Range("BA3").Value = 51
$ws.Range("BB3").Value = 67
$ws.Range("BC3").Value = 151
